$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds a "Förändrad" (last changed) date for each record.
# It was bumped by one day, from 2023-09-12 (serial 45181) to 2023-09-13 (serial 45182),
# for every data row (rows 2 through 452).
$ws.Range("C2:C452").Value = 45182
